$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy formatting from the existing G1
# header cell (bold, bordered, centered) then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H5 with 0 (plain numeric cells, same as existing data columns)
$ws.Range("H2:H5").Value = 0
